# Fix a missing period ("corrigindo erros ortográficos na documentação do
# projeto") at the end of the paragraph about barrel aging / wine storage.
# The paragraph currently ends in "...posicionados horizontalmente" with no
# final punctuation; we append a "." styled the same as the rest of the
# paragraph (Arial font, black/text1 color) as a new run right after the
# existing text run.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(
    "posicionados horizontalmente",  # Find what
    $true,                            # MatchCase
    $false,                           # MatchWholeWord
    $false,                           # MatchWildcards
    $false,                           # MatchSoundsLike
    $false,                           # MatchAllWordForms
    $true,                            # Forward
    1,                                 # Wrap (wdFindContinue)
    $false,                           # Format
    "",                                # Replace with
    0                                  # Replace (wdReplaceNone)
)

if (-not $found) {
    throw "Could not find target sentence to fix punctuation."
}

# Collapse to the end of the found text and insert the missing period as a
# new run, then explicitly (re)apply the paragraph's character formatting
# (Arial / Arial for complex-script + the black "text1" theme color) so the
# appended run keeps the same visual style as the sentence it completes.
$rng.Collapse(0)
$rng.InsertAfter(".")
$rng.Font.Name = "Arial"
$rng.Font.NameBi = "Arial"
